# Update the "About" sheet's currency-year text from 2021 to 2023 in all places.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# A18: "billion 2021 dollars" -> "billion 2023 dollars"
$ws.Range("A18").Value = "billion 2023 dollars"

# A21: "million 2021 dollars" -> "million 2023 dollars"
$ws.Range("A21").Value = "million 2023 dollars"

# A24: "2021 dollars" -> "2023 dollars"
$ws.Range("A24").Value = "2023 dollars"

# B29: "which in this case is "2012 dollars per 2021 dollar."" -> "...2023 dollar.""
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2023 dollar."'

# Selection moved to B30 in the saved file
$ws.Range("B30").Select()
